# Website Testing.xlsx — update compatibility-testing comment text and
# move the active selection (Testing, Classification Results and Corpus).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Compability Testing" section (rows 24-28) all share the same
# comment text in column C. Update each cell to the new wording — Excel
# dedupes identical strings in the shared-string table, so this updates
# the single shared entry used by all five cells.
$newComment = "Tested in Chrome, Internet Explorer, Safari and Firefox"
$ws.Range("C24").Value = $newComment
$ws.Range("C25").Value = $newComment
$ws.Range("C26").Value = $newComment
$ws.Range("C27").Value = $newComment
$ws.Range("C28").Value = $newComment

# Move the selection/active cell down to C30 (also clears the scrolled
# topLeftCell position that was previously pinned at A20).
$ws.Range("C30").Select()
